$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Fix constraints on player when game paused") {
        $p.Range.Delete()
    }
}
